$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("E2").Value = 0.03646020042922023
$ws.Range("F2").Value = 0.1874935731250425
$ws.Range("G2").Value = 0.2663658411335913
$ws.Range("H2").Value = 0.2838937758016715
$ws.Range("L2").Value = 0.1128397561067585
$ws.Range("M2").Value = 0.02490907487844918
$ws.Range("N2").Value = 0.01017883168603257
$ws.Range("Q2").Value = 0.0003408374543019545
$ws.Range("R2").Value = 0.005539469710326634
$ws.Range("S2").Value = 0.006668552896578233
$ws.Range("T2").Value = 0.002591626713728028
$ws.Range("X2").Value = 0.01759711984284027
$ws.Range("Y2").Value = 0.01074089004420585
$ws.Range("AC2").Value = 0.03438045017725342
$ws.Range("E3").Value = 0.2908232273162981
$ws.Range("G3").Value = 0.3586975870839298
$ws.Range("H3").Value = 0.03855900497585855
$ws.Range("J3").Value = 0.01858736019769114
$ws.Range("K3").Value = 0.1208489869191962
$ws.Range("L3").Value = 0.03914059201849913
$ws.Range("M3").Value = 0.02201482564134247
$ws.Range("N3").Value = 0.01054747066721586
$ws.Range("O3").Value = 0.0125846300953745
$ws.Range("R3").Value = 0.01957770292497812
$ws.Range("X3").Value = 0.03610085126994152
$ws.Range("AA3").Value = 0.00129372728321139
$ws.Range("AB3").Value = 0.005792561281552891
$ws.Range("AC3").Value = 0.02543147232491038
$ws.Range("E4").Value = 0.3115819456258839
$ws.Range("F4").Value = 0.05839872578124777
$ws.Range("G4").Value = 0.3360468485982868
$ws.Range("H4").Value = 0.009140668295738685
$ws.Range("J4").Value = 0.002067701341536079
$ws.Range("K4").Value = 0.1552259915244661
$ws.Range("L4").Value = 0.03470525542210808
$ws.Range("M4").Value = 0.03701417222022348
$ws.Range("O4").Value = 0.0011142937367423
$ws.Range("R4").Value = 0.009279093643155368
$ws.Range("X4").Value = 0.02465399748063748
$ws.Range("AB4").Value = 0.01223352148028301
$ws.Range("AC4").Value = 0.002759091809743823
$ws.Range("AD4").Value = 0.005778693039947349
$ws.Range("E5").Value = 0.184375685329555
$ws.Range("G5").Value = 0.3879745534065988
$ws.Range("H5").Value = 0.111850781774743
$ws.Range("J5").Value = 0.02170350158576591
$ws.Range("K5").Value = 0.04567886471806309
$ws.Range("L5").Value = 0.1203816505237333
$ws.Range("M5").Value = 0.05030297971788403
$ws.Range("N5").Value = 0.002593960872364038
$ws.Range("P5").Value = 0.0001684614241921933
$ws.Range("R5").Value = 0.005089454327890633
$ws.Range("T5").Value = 0.003235602090665387
$ws.Range("X5").Value = 0.04431710709177726
$ws.Range("AA5").Value = 0.003901752111189831
$ws.Range("AC5").Value = 0.01842564502557766
$ws.Range("D6").Value = 0.0146589640158236
$ws.Range("E6").Value = 0.2314252557623172
$ws.Range("F6").Value = 0.2524613736564633
$ws.Range("G6").Value = 0.2021378390062248
$ws.Range("I6").Value = 0.004924783326422265
$ws.Range("K6").Value = 0.1593934997664758
$ws.Range("L6").Value = 0.0283228084378554
$ws.Range("M6").Value = 0.04792730531266454
$ws.Range("P6").Value = 0.01070484535079365
$ws.Range("W6").Value = 0.01448531176682427
$ws.Range("X6").Value = 0.01493110183893391
$ws.Range("AB6").Value = 0.01862691175920119

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("E2").Value = 0.03646020042922023
$ws.Range("F2").Value = 0.2239537735542627
$ws.Range("G2").Value = 0.490319614687854
$ws.Range("H2").Value = 0.7742133904895255
$ws.Range("I2").Value = 0.7742133904895255
$ws.Range("J2").Value = 0.7742133904895255
$ws.Range("K2").Value = 0.7742133904895255
$ws.Range("L2").Value = 0.887053146596284
$ws.Range("M2").Value = 0.9119622214747332
$ws.Range("N2").Value = 0.9221410531607658
$ws.Range("O2").Value = 0.9221410531607658
$ws.Range("P2").Value = 0.9221410531607658
$ws.Range("Q2").Value = 0.9224818906150677
$ws.Range("R2").Value = 0.9280213603253943
$ws.Range("S2").Value = 0.9346899132219725
$ws.Range("T2").Value = 0.9372815399357005
$ws.Range("U2").Value = 0.9372815399357005
$ws.Range("V2").Value = 0.9372815399357005
$ws.Range("W2").Value = 0.9372815399357005
$ws.Range("X2").Value = 0.9548786597785408
$ws.Range("Y2").Value = 0.9656195498227467
$ws.Range("Z2").Value = 0.9656195498227467
$ws.Range("AA2").Value = 0.9656195498227467
$ws.Range("AB2").Value = 0.9656195498227467
$ws.Range("E3").Value = 0.2908232273162981
$ws.Range("F3").Value = 0.2908232273162981
$ws.Range("G3").Value = 0.6495208144002278
$ws.Range("H3").Value = 0.6880798193760863
$ws.Range("I3").Value = 0.6880798193760863
$ws.Range("J3").Value = 0.7066671795737774
$ws.Range("K3").Value = 0.8275161664929737
$ws.Range("L3").Value = 0.8666567585114728
$ws.Range("M3").Value = 0.8886715841528153
$ws.Range("N3").Value = 0.8992190548200312
$ws.Range("O3").Value = 0.9118036849154056
$ws.Range("P3").Value = 0.9118036849154056
$ws.Range("Q3").Value = 0.9118036849154056
$ws.Range("R3").Value = 0.9313813878403837
$ws.Range("S3").Value = 0.9313813878403837
$ws.Range("T3").Value = 0.9313813878403837
$ws.Range("U3").Value = 0.9313813878403837
$ws.Range("V3").Value = 0.9313813878403837
$ws.Range("W3").Value = 0.9313813878403837
$ws.Range("X3").Value = 0.9674822391103253
$ws.Range("Y3").Value = 0.9674822391103253
$ws.Range("Z3").Value = 0.9674822391103253
$ws.Range("AA3").Value = 0.9687759663935367
$ws.Range("AB3").Value = 0.9745685276750896
$ws.Range("E4").Value = 0.3115819456258839
$ws.Range("F4").Value = 0.3699806714071317
$ws.Range("G4").Value = 0.7060275200054185
$ws.Range("H4").Value = 0.7151681883011571
$ws.Range("I4").Value = 0.7151681883011571
$ws.Range("J4").Value = 0.7172358896426932
$ws.Range("K4").Value = 0.8724618811671594
$ws.Range("L4").Value = 0.9071671365892675
$ws.Range("M4").Value = 0.944181308809491
$ws.Range("N4").Value = 0.944181308809491
$ws.Range("O4").Value = 0.9452956025462333
$ws.Range("P4").Value = 0.9452956025462333
$ws.Range("Q4").Value = 0.9452956025462333
$ws.Range("R4").Value = 0.9545746961893887
$ws.Range("S4").Value = 0.9545746961893887
$ws.Range("T4").Value = 0.9545746961893887
$ws.Range("U4").Value = 0.9545746961893887
$ws.Range("V4").Value = 0.9545746961893887
$ws.Range("W4").Value = 0.9545746961893887
$ws.Range("X4").Value = 0.9792286936700262
$ws.Range("Y4").Value = 0.9792286936700262
$ws.Range("Z4").Value = 0.9792286936700262
$ws.Range("AA4").Value = 0.9792286936700262
$ws.Range("AB4").Value = 0.9914622151503092
$ws.Range("AC4").Value = 0.9942213069600531
$ws.Range("E5").Value = 0.184375685329555
$ws.Range("F5").Value = 0.184375685329555
$ws.Range("G5").Value = 0.5723502387361536
$ws.Range("H5").Value = 0.6842010205108966
$ws.Range("I5").Value = 0.6842010205108966
$ws.Range("J5").Value = 0.7059045220966625
$ws.Range("K5").Value = 0.7515833868147256
$ws.Range("L5").Value = 0.8719650373384589
$ws.Range("M5").Value = 0.9222680170563429
$ws.Range("N5").Value = 0.9248619779287069
$ws.Range("O5").Value = 0.9248619779287069
$ws.Range("P5").Value = 0.9250304393528992
$ws.Range("Q5").Value = 0.9250304393528992
$ws.Range("R5").Value = 0.9301198936807898
$ws.Range("S5").Value = 0.9301198936807898
$ws.Range("T5").Value = 0.9333554957714553
$ws.Range("U5").Value = 0.9333554957714553
$ws.Range("V5").Value = 0.9333554957714553
$ws.Range("W5").Value = 0.9333554957714553
$ws.Range("X5").Value = 0.9776726028632325
$ws.Range("Y5").Value = 0.9776726028632325
$ws.Range("Z5").Value = 0.9776726028632325
$ws.Range("AA5").Value = 0.9815743549744224
$ws.Range("AB5").Value = 0.9815743549744224
$ws.Range("D6").Value = 0.0146589640158236
$ws.Range("E6").Value = 0.2460842197781408
$ws.Range("F6").Value = 0.4985455934346041
$ws.Range("G6").Value = 0.700683432440829
$ws.Range("H6").Value = 0.700683432440829
$ws.Range("I6").Value = 0.7056082157672512
$ws.Range("J6").Value = 0.7056082157672512
$ws.Range("K6").Value = 0.8650017155337271
$ws.Range("L6").Value = 0.8933245239715825
$ws.Range("M6").Value = 0.941251829284247
$ws.Range("N6").Value = 0.941251829284247
$ws.Range("O6").Value = 0.941251829284247
$ws.Range("P6").Value = 0.9519566746350406
$ws.Range("Q6").Value = 0.9519566746350406
$ws.Range("R6").Value = 0.9519566746350406
$ws.Range("S6").Value = 0.9519566746350406
$ws.Range("T6").Value = 0.9519566746350406
$ws.Range("U6").Value = 0.9519566746350406
$ws.Range("V6").Value = 0.9519566746350406
$ws.Range("W6").Value = 0.966441986401865
$ws.Range("X6").Value = 0.9813730882407988
$ws.Range("Y6").Value = 0.9813730882407988
$ws.Range("Z6").Value = 0.9813730882407988
$ws.Range("AA6").Value = 0.9813730882407988
$ws.Range("AB6").Value = 1
$ws.Range("AC6").Value = 1
$ws.Range("AD6").Value = 1
$ws.Range("AE6").Value = 1
$ws.Range("AF6").Value = 1
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 1
$ws.Range("AJ6").Value = 1

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 0.7742133904895255
$ws.Range("G2").Value = 4
$ws.Range("F3").Value = 0.6495208144002278
$ws.Range("F4").Value = 0.7060275200054185
$ws.Range("F5").Value = 0.5723502387361536
$ws.Range("D6").Value = 6
$ws.Range("F6").Value = 0.700683432440829
$ws.Range("G6").Value = 4

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.7742133904895255
$ws.Range("D3").Value = 9
$ws.Range("F3").Value = 0.7066671795737774
$ws.Range("G3").Value = 7
$ws.Range("F4").Value = 0.7060275200054185
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.7059045220966625
$ws.Range("G5").Value = 7
$ws.Range("F6").Value = 0.700683432440829

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 11
$ws.Range("F2").Value = 0.887053146596284
$ws.Range("G2").Value = 8
$ws.Range("D3").Value = 10
$ws.Range("F3").Value = 0.8275161664929737
$ws.Range("G3").Value = 8
$ws.Range("D4").Value = 10
$ws.Range("F4").Value = 0.8724618811671594
$ws.Range("G4").Value = 8
$ws.Range("D5").Value = 11
$ws.Range("F5").Value = 0.8719650373384589
$ws.Range("G5").Value = 9
$ws.Range("F6").Value = 0.8650017155337271

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 12
$ws.Range("F2").Value = 0.9119622214747332
$ws.Range("G2").Value = 9
$ws.Range("D3").Value = 14
$ws.Range("F3").Value = 0.9118036849154056
$ws.Range("G3").Value = 12
$ws.Range("D4").Value = 11
$ws.Range("F4").Value = 0.9071671365892675
$ws.Range("G4").Value = 9
$ws.Range("D5").Value = 12
$ws.Range("F5").Value = 0.9222680170563429
$ws.Range("G5").Value = 10
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.941251829284247
$ws.Range("G6").Value = 10
